$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04868399285477
$ws.Range("D2").Value = 1.054860339134333
$ws.Range("E2").Value = 1.052356275818464
$ws.Range("F2").Value = 1.063517475306389
$ws.Range("I2").Value = 1.04425969776581
$ws.Range("J2").Value = 1.053726439327648
$ws.Range("K2").Value = 1.057602386701298
$ws.Range("L2").Value = 1.055105234753143
$ws.Range("M2").Value = 1.066235900098111
$ws.Range("N2").Value = 1.021481688717232

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050179542769677
$ws.Range("D3").Value = 1.056048379611261
$ws.Range("E3").Value = 1.053802302063283
$ws.Range("F3").Value = 1.064841305221988
$ws.Range("I3").Value = 1.044680743021511
$ws.Range("J3").Value = 1.054868289139907
$ws.Range("K3").Value = 1.058602797307515
$ws.Range("L3").Value = 1.056362465725956
$ws.Range("M3").Value = 1.067373479898062
$ws.Range("N3").Value = 1.021884764034187

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.05114574237534
$ws.Range("D4").Value = 1.056815616780974
$ws.Range("E4").Value = 1.054736822203036
$ws.Range("F4").Value = 1.065696655898541
$ws.Range("I4").Value = 1.044950900949727
$ws.Range("J4").Value = 1.055605225164153
$ws.Range("K4").Value = 1.059248056116493
$ws.Range("L4").Value = 1.057174308922963
$ws.Range("M4").Value = 1.068107769947436
$ws.Range("N4").Value = 1.022144351324956

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051551576218142
$ws.Range("D5").Value = 1.057137808254773
$ws.Range("E5").Value = 1.055129424573855
$ws.Range("F5").Value = 1.066055950338016
$ws.Range("I5").Value = 1.045063930312082
$ws.Range("J5").Value = 1.055914579800928
$ws.Range("K5").Value = 1.059518831271309
$ws.Range("L5").Value = 1.057515214287037
$ws.Range("M5").Value = 1.06841603973656
$ws.Range("N5").Value = 1.022253189505872

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051619696825591
$ws.Range("D6").Value = 1.057191884950341
$ws.Range("E6").Value = 1.055195328607812
$ws.Range("F6").Value = 1.066116260321004
$ws.Range("I6").Value = 1.04508287654684
$ws.Range("J6").Value = 1.055966495390931
$ws.Range("K6").Value = 1.059564266919963
$ws.Range("L6").Value = 1.057572430898948
$ws.Range("M6").Value = 1.068467774775359
$ws.Range("N6").Value = 1.022271446831811

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051151166538162
$ws.Range("D7").Value = 1.056819923307907
$ws.Range("E7").Value = 1.054742069227292
$ws.Range("F7").Value = 1.065701457960526
$ws.Range("I7").Value = 1.044952413392293
$ws.Range("J7").Value = 1.055609360549937
$ws.Range("K7").Value = 1.059251676153275
$ws.Range("L7").Value = 1.057178865655805
$ws.Range("M7").Value = 1.068111890729911
$ws.Range("N7").Value = 1.022145806772173

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049189740573841
$ws.Range("D8").Value = 1.05526215758769
$ws.Range("E8").Value = 1.052845210487636
$ws.Range("F8").Value = 1.06396513326168
$ws.Range("I8").Value = 1.044402466845946
$ws.Range("J8").Value = 1.054112733220881
$ws.Range("K8").Value = 1.057940912438417
$ws.Range("L8").Value = 1.05553047118295
$ws.Range("M8").Value = 1.066620726397469
$ws.Range("N8").Value = 1.021618165670216

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045721460649754
$ws.Range("D9").Value = 1.052505421546323
$ws.Range("E9").Value = 1.049493549956164
$ws.Range("F9").Value = 1.06089563472792
$ws.Range("I9").Value = 1.04341577218937
$ws.Range("J9").Value = 1.051460544967241
$ws.Range("K9").Value = 1.055615080187569
$ws.Range("L9").Value = 1.052612729341264
$ws.Range("M9").Value = 1.063979082663324
$ws.Range("N9").Value = 1.020678888888288

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043400721997567
$ws.Range("D10").Value = 1.050659368637467
$ws.Range("E10").Value = 1.047252529611018
$ws.Range("F10").Value = 1.058842305877132
$ws.Range("I10").Value = 1.042745979172082
$ws.Range("J10").Value = 1.049682014033364
$ws.Range("K10").Value = 1.054053395849909
$ws.Range("L10").Value = 1.050658400170753
$ws.Range("M10").Value = 1.062208228803131
$ws.Range("N10").Value = 1.020046186412452

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042393679377517
$ws.Range("D11").Value = 1.049857982163744
$ws.Range("E11").Value = 1.046280486933973
$ws.Range("F11").Value = 1.057951451713773
$ws.Range("I11").Value = 1.04245307188364
$ws.Range("J11").Value = 1.048909341547458
$ws.Range("K11").Value = 1.053374460847335
$ws.Range("L11").Value = 1.049809892786797
$ws.Range("M11").Value = 1.061439040714782
$ws.Range("N11").Value = 1.019770644783733

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042019286831054
$ws.Range("D12").Value = 1.0495599996319
$ws.Range("E12").Value = 1.045919168389354
$ws.Range("F12").Value = 1.057620278969486
$ws.Range("I12").Value = 1.042343836999897
$ws.Range("J12").Value = 1.0486219453852
$ws.Range("K12").Value = 1.053121860191005
$ws.Range("L12").Value = 1.049494370768461
$ws.Range("M12").Value = 1.061152963434989
$ws.Range("N12").Value = 1.019668056886031

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042099610482867
$ws.Range("D13").Value = 1.04962393211495
$ws.Range("E13").Value = 1.045996684280565
$ws.Range("F13").Value = 1.057691328971981
$ws.Range("I13").Value = 1.042367288035258
$ws.Range("J13").Value = 1.048683610654053
$ws.Range("K13").Value = 1.053176062719201
$ws.Range("L13").Value = 1.049562067209198
$ws.Range("M13").Value = 1.061214344675912
$ws.Range("N13").Value = 1.019690073215107

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042362738802037
$ws.Range("D14").Value = 1.049833357234994
$ws.Range("E14").Value = 1.046250625554281
$ws.Range("F14").Value = 1.057924082417078
$ws.Range("I14").Value = 1.042444051407404
$ws.Range("J14").Value = 1.048885593312383
$ws.Range("K14").Value = 1.053353589278041
$ws.Range("L14").Value = 1.049783818793538
$ws.Range("M14").Value = 1.06141540099343
$ws.Range("N14").Value = 1.019762169735712

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042524816556346
$ws.Range("D15").Value = 1.049962349438713
$ws.Range("E15").Value = 1.046407052635709
$ws.Range("F15").Value = 1.058067453527973
$ws.Range("I15").Value = 1.042491289979688
$ws.Range("J15").Value = 1.049009989444816
$ws.Range("K15").Value = 1.053462914214715
$ws.Range("L15").Value = 1.049920400843802
$ws.Range("M15").Value = 1.061539229656077
$ws.Range("N15").Value = 1.019806558974579

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043467510050314
$ws.Range("D16").Value = 1.050712510630593
$ws.Range("E16").Value = 1.047317005072408
$ws.Range("F16").Value = 1.058901391480197
$ws.Range("I16").Value = 1.042765357476789
$ws.Range("J16").Value = 1.049733239290255
$ws.Range("K16").Value = 1.054098396808129
$ws.Range("L16").Value = 1.050714664341956
$ws.Range("M16").Value = 1.062259226234304
$ws.Range("N16").Value = 1.020064439758431

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044058255188459
$ws.Range("D17").Value = 1.051182518017036
$ws.Range("E17").Value = 1.047887342295652
$ws.Range("F17").Value = 1.059424025455829
$ws.Range("I17").Value = 1.042936498933227
$ws.Range("J17").Value = 1.050186225301155
$ws.Range("K17").Value = 1.0544962870373
$ws.Range("L17").Value = 1.051212271798858
$ws.Range("M17").Value = 1.062710214833255
$ws.Range("N17").Value = 1.020225777518911

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044402620275487
$ws.Range("D18").Value = 1.051456469619062
$ws.Range("E18").Value = 1.04821985016103
$ws.Range("F18").Value = 1.059728700716539
$ws.Range("I18").Value = 1.043036044935979
$ws.Range("J18").Value = 1.050450198158207
$ws.Range("K18").Value = 1.054728108090853
$ws.Range("L18").Value = 1.051502299306063
$ws.Range("M18").Value = 1.062973038050587
$ws.Range("N18").Value = 1.02031973109023

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044520005119694
$ws.Range("D19").Value = 1.051549847025555
$ws.Range("E19").Value = 1.048333199829743
$ws.Range("F19").Value = 1.059832558773008
$ws.Range("I19").Value = 1.043069940502075
$ws.Range("J19").Value = 1.050540164530217
$ws.Range("K19").Value = 1.054807108890749
$ws.Range("L19").Value = 1.051601154356627
$ws.Range("M19").Value = 1.063062615023233
$ws.Range("N19").Value = 1.020351741120963

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043994895259537
$ws.Range("D20").Value = 1.05113211096695
$ws.Range("E20").Value = 1.047826167138366
$ws.Range("F20").Value = 1.059367969235006
$ws.Range("I20").Value = 1.042918165842342
$ws.Range("J20").Value = 1.050137649707974
$ws.Range("K20").Value = 1.054453624281887
$ws.Range("L20").Value = 1.051158905883693
$ws.Range("M20").Value = 1.062661851944803
$ws.Range("N20").Value = 1.02020848324583

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042285263313558
$ws.Range("D21").Value = 1.049771695419117
$ws.Range("E21").Value = 1.046175853384525
$ws.Range("F21").Value = 1.057855549834554
$ws.Range("I21").Value = 1.04242145856973
$ws.Range("J21").Value = 1.048826125309733
$ws.Range("K21").Value = 1.053301323607241
$ws.Range("L21").Value = 1.049718528162939
$ws.Range("M21").Value = 1.061356205090888
$ws.Range("N21").Value = 1.019740945738015

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.041208423817898
$ws.Range("D22").Value = 1.048914540311687
$ws.Range("E22").Value = 1.045136735204493
$ws.Range("F22").Value = 1.05690306697166
$ws.Range("I22").Value = 1.042106634484291
$ws.Range("J22").Value = 1.047999250979627
$ws.Range("K22").Value = 1.052574428339141
$ws.Range("L22").Value = 1.048810884863499
$ws.Range("M22").Value = 1.060533168938667
$ws.Range("N22").Value = 1.019445599609166

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.041779462546431
$ws.Range("D23").Value = 1.04936910797362
$ws.Range("E23").Value = 1.045687736423863
$ws.Range("F23").Value = 1.057408146845643
$ws.Range("I23").Value = 1.042273768995838
$ws.Range("J23").Value = 1.048437809952343
$ws.Range("K23").Value = 1.052959998606046
$ws.Range("L23").Value = 1.049292237823789
$ws.Range("M23").Value = 1.060969679504866
$ws.Range("N23").Value = 1.019602300466257

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044023525522575
$ws.Range("D24").Value = 1.051154888348266
$ws.Range("E24").Value = 1.047853810053906
$ws.Range("F24").Value = 1.059393299148583
$ws.Range("I24").Value = 1.042926450635821
$ws.Range("J24").Value = 1.050159599688623
$ws.Range("K24").Value = 1.054472902553151
$ws.Range("L24").Value = 1.051183020319028
$ws.Range("M24").Value = 1.062683705765739
$ws.Range("N24").Value = 1.020216298253223

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04661956565824
$ws.Range("D25").Value = 1.053219530201296
$ws.Range("E25").Value = 1.050361162894234
$ws.Range("F25").Value = 1.061690380378365
$ws.Range("I25").Value = 1.043672959743015
$ws.Range("J25").Value = 1.052148006614294
$ws.Range("K25").Value = 1.056218302683831
$ws.Range("L25").Value = 1.053368622910957
$ws.Range("M25").Value = 1.064663706402571
$ws.Range("N25").Value = 1.020922853933943
